$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Fitness" values for rows 2..74 (column C), per the target diff.
$values = @(
    11683, 10943, 10943, 10943, 10901, 10901, 10517, 10517, 10517, 10517,
    10517, 10517, 10495, 10053, 10053, 10053, 9990, 9990, 9990, 9990,
    9990, 9990, 9990, 9990, 9990, 9687, 9687, 9687, 9547, 9547,
    9547, 9547, 9110, 9110, 8789, 8789, 8699, 8699, 8694, 8694,
    8694, 8404, 8404, 8404, 8404, 8404, 8404, 8404, 8404, 8404,
    8404, 8404, 7920, 7920, 7920, 7920, 7920, 7682, 7682, 7682,
    7682, 7682, 7682, 7682, 7682, 7674, 7674, 7674, 7657, 7657,
    7657, 7657, 7657
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
